# Applies the updated statistics (rows 4-13) to Sheet1.
# Each row block lists (column-index, new-value) pairs derived from the
# authoritative cell-level diff between the old and new workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row4 = @(
    @(2, 0.307),  # B4
    @(3, 0.051),  # C4
    @(4, 0.227),  # D4
    @(5, 0.16),  # E4
    @(7, 0.112),  # G4
    @(8, 0.2),  # H4
    @(10, 0.101),  # J4
    @(11, 0.361),  # K4
    @(12, 0.101),  # L4
    @(13, 0.318),  # M4
    @(14, 0.28),  # N4
    @(15, 0.019),  # O4
    @(16, 0.139),  # P4
    @(17, 0.542),  # Q4
    @(18, 0.214),  # R4
    @(19, 0.462),  # S4
    @(20, 0.299),  # T4
    @(23, 0.258),  # W4
    @(24, 0.042),  # X4
    @(25, 0.205),  # Y4
    @(26, 0.476),  # Z4
    @(27, 0.128),  # AA4
    @(28, 0.358),  # AB4
    @(31, 0.076),  # AE4
    @(32, 0.743),  # AF4
    @(34, 0.301),  # AH4
    @(35, 0.681),  # AI4
    @(36, 0.169),  # AJ4
    @(37, 0.411),  # AK4
    @(38, 0.732),  # AL4
    @(39, 0.103),  # AM4
    @(40, 0.321),  # AN4
    @(41, 0.719)  # AO4
)

$row5 = @(
    @(2, 0.833),  # B5
    @(3, 0.139),  # C5
    @(4, 0.373),  # D5
    @(5, 0.694),  # E5
    @(6, 0.212),  # F5
    @(7, 0.461),  # G5
    @(8, 0.861),  # H5
    @(9, 0.12),  # I5
    @(10, 0.346),  # J5
    @(11, 0.694),  # K5
    @(12, 0.212),  # L5
    @(13, 0.461),  # M5
    @(14, 0.861),  # N5
    @(15, 0.12),  # O5
    @(16, 0.346),  # P5
    @(17, 0.611),  # Q5
    @(18, 0.238),  # R5
    @(19, 0.487),  # S5
    @(20, 0.611),  # T5
    @(21, 0.238),  # U5
    @(22, 0.487),  # V5
    @(23, 0.778),  # W5
    @(24, 0.173),  # X5
    @(25, 0.416),  # Y5
    @(26, 0.861),  # Z5
    @(27, 0.12),  # AA5
    @(28, 0.346),  # AB5
    @(29, 0.778),  # AC5
    @(30, 0.173),  # AD5
    @(31, 0.416),  # AE5
    @(32, 0.972),  # AF5
    @(33, 0.027),  # AG5
    @(34, 0.164),  # AH5
    @(35, 0.778),  # AI5
    @(36, 0.173),  # AJ5
    @(37, 0.416),  # AK5
    @(38, 0.944),  # AL5
    @(39, 0.052),  # AM5
    @(40, 0.229),  # AN5
    @(41, 0.898)  # AO5
)

$row6 = @(
    @(2, 0.449),  # B6
    @(5, 0.26),  # E6
    @(8, 0.325),  # H6
    @(11, 0.475),  # K6
    @(14, 0.423),  # N6
    @(17, 0.574),  # Q6
    @(20, 0.402),  # T6
    @(23, 0.387),  # W6
    @(26, 0.613),  # Z6
    @(32, 0.842),  # AF6
    @(35, 0.726),  # AI6
    @(38, 0.825),  # AL6
    @(41, 0.798)  # AO6
)

$row7 = @(
    @(2, 0.62),  # B7
    @(5, 0.416),  # E7
    @(8, 0.518),  # H7
    @(11, 0.586),  # K7
    @(14, 0.608),  # N7
    @(17, 0.596),  # Q7
    @(20, 0.506),  # T7
    @(23, 0.554),  # W7
    @(26, 0.741),  # Z7
    @(29, 0.384),  # AC7
    @(32, 0.916),  # AF7
    @(35, 0.756),  # AI7
    @(38, 0.892),  # AL7
    @(41, 0.855)  # AO7
)

$row8 = @(
    @(2, 0.765),  # B8
    @(3, 0.14),  # C8
    @(4, 0.375),  # D8
    @(5, 0.578),  # E8
    @(8, 0.742),  # H8
    @(9, 0.131),  # I8
    @(10, 0.363),  # J8
    @(11, 0.618),  # K8
    @(12, 0.195),  # L8
    @(13, 0.442),  # M8
    @(14, 0.777),  # N8
    @(15, 0.127),  # O8
    @(16, 0.356),  # P8
    @(17, 0.58),  # Q8
    @(18, 0.224),  # R8
    @(19, 0.473),  # S8
    @(20, 0.528),  # T8
    @(21, 0.202),  # U8
    @(22, 0.45),  # V8
    @(23, 0.699),  # W8
    @(24, 0.164),  # X8
    @(25, 0.405),  # Y8
    @(26, 0.792),  # Z8
    @(27, 0.125),  # AA8
    @(28, 0.353),  # AB8
    @(29, 0.663),  # AC8
    @(30, 0.171),  # AD8
    @(31, 0.414),  # AE8
    @(32, 0.897),  # AF8
    @(33, 0.047),  # AG8
    @(34, 0.216),  # AH8
    @(35, 0.768),  # AI8
    @(36, 0.172),  # AJ8
    @(37, 0.415),  # AK8
    @(38, 0.914),  # AL8
    @(39, 0.059),  # AM8
    @(40, 0.244),  # AN8
    @(41, 0.86)  # AO8
)

$row9 = @(
    @(2, 0.667),  # B9
    @(3, 0.222),  # C9
    @(4, 0.471),  # D9
    @(5, 0.444),  # E9
    @(6, 0.247),  # F9
    @(7, 0.497),  # G9
    @(8, 0.611),  # H9
    @(9, 0.238),  # I9
    @(10, 0.487),  # J9
    @(11, 0.528),  # K9
    @(12, 0.249),  # L9
    @(13, 0.499),  # M9
    @(14, 0.667),  # N9
    @(15, 0.222),  # O9
    @(16, 0.471),  # P9
    @(17, 0.528),  # Q9
    @(18, 0.249),  # R9
    @(19, 0.499),  # S9
    @(20, 0.417),  # T9
    @(21, 0.243),  # U9
    @(22, 0.493),  # V9
    @(23, 0.583),  # W9
    @(24, 0.243),  # X9
    @(25, 0.493),  # Y9
    @(26, 0.694),  # Z9
    @(27, 0.212),  # AA9
    @(28, 0.461),  # AB9
    @(29, 0.556),  # AC9
    @(30, 0.247),  # AD9
    @(31, 0.497),  # AE9
    @(32, 0.778),  # AF9
    @(33, 0.173),  # AG9
    @(34, 0.416),  # AH9
    @(35, 0.75),  # AI9
    @(36, 0.188),  # AJ9
    @(37, 0.433),  # AK9
    @(38, 0.861),  # AL9
    @(39, 0.12),  # AM9
    @(40, 0.346),  # AN9
    @(41, 0.796)  # AO9
)

$row10 = @(
    @(2, 0.833),  # B10
    @(3, 0.139),  # C10
    @(4, 0.373),  # D10
    @(5, 0.611),  # E10
    @(6, 0.238),  # F10
    @(7, 0.487),  # G10
    @(8, 0.778),  # H10
    @(9, 0.173),  # I10
    @(10, 0.416),  # J10
    @(11, 0.694),  # K10
    @(12, 0.212),  # L10
    @(13, 0.461),  # M10
    @(14, 0.833),  # N10
    @(15, 0.139),  # O10
    @(16, 0.373),  # P10
    @(17, 0.611),  # Q10
    @(18, 0.238),  # R10
    @(19, 0.487),  # S10
    @(20, 0.611),  # T10
    @(21, 0.238),  # U10
    @(22, 0.487),  # V10
    @(23, 0.778),  # W10
    @(24, 0.173),  # X10
    @(25, 0.416),  # Y10
    @(26, 0.861),  # Z10
    @(27, 0.12),  # AA10
    @(28, 0.346),  # AB10
    @(29, 0.667),  # AC10
    @(30, 0.222),  # AD10
    @(31, 0.471),  # AE10
    @(32, 0.972),  # AF10
    @(33, 0.027),  # AG10
    @(34, 0.164),  # AH10
    @(35, 0.778),  # AI10
    @(36, 0.173),  # AJ10
    @(37, 0.416),  # AK10
    @(38, 0.944),  # AL10
    @(39, 0.052),  # AM10
    @(40, 0.229),  # AN10
    @(41, 0.898)  # AO10
)

$row11 = @(
    @(2, 0.833),  # B11
    @(3, 0.139),  # C11
    @(4, 0.373),  # D11
    @(5, 0.694),  # E11
    @(6, 0.212),  # F11
    @(7, 0.461),  # G11
    @(8, 0.861),  # H11
    @(9, 0.12),  # I11
    @(10, 0.346),  # J11
    @(11, 0.694),  # K11
    @(12, 0.212),  # L11
    @(13, 0.461),  # M11
    @(14, 0.861),  # N11
    @(15, 0.12),  # O11
    @(16, 0.346),  # P11
    @(17, 0.611),  # Q11
    @(18, 0.238),  # R11
    @(19, 0.487),  # S11
    @(20, 0.611),  # T11
    @(21, 0.238),  # U11
    @(22, 0.487),  # V11
    @(23, 0.778),  # W11
    @(24, 0.173),  # X11
    @(25, 0.416),  # Y11
    @(26, 0.861),  # Z11
    @(27, 0.12),  # AA11
    @(28, 0.346),  # AB11
    @(29, 0.722),  # AC11
    @(30, 0.201),  # AD11
    @(31, 0.448),  # AE11
    @(32, 0.972),  # AF11
    @(33, 0.027),  # AG11
    @(34, 0.164),  # AH11
    @(35, 0.778),  # AI11
    @(36, 0.173),  # AJ11
    @(37, 0.416),  # AK11
    @(38, 0.944),  # AL11
    @(39, 0.052),  # AM11
    @(40, 0.229),  # AN11
    @(41, 0.898)  # AO11
)

$row12 = @(
    @(2, 1.267),  # B12
    @(3, 0.329),  # C12
    @(4, 0.573),  # D12
    @(5, 1.68),  # E12
    @(6, 1.098),  # F12
    @(7, 1.048),  # G12
    @(8, 1.613),  # H12
    @(9, 1.334),  # I12
    @(10, 1.155),  # J12
    @(11, 1.4),  # K12
    @(12, 0.5600000000000001),  # L12
    @(13, 0.748),  # M12
    @(14, 1.355),  # N12
    @(15, 0.552),  # O12
    @(16, 0.743),  # P12
    @(26, 1.258),  # Z12
    @(27, 0.32),  # AA12
    @(28, 0.5659999999999999),  # AB12
    @(29, 1.821),  # AC12
    @(30, 2.504),  # AD12
    @(31, 1.582),  # AE12
    @(32, 1.229),  # AF12
    @(33, 0.233),  # AG12
    @(34, 0.483),  # AH12
    @(35, 1.036),  # AI12
    @(36, 0.034),  # AJ12
    @(37, 0.186),  # AK12
    @(38, 1.088),  # AL12
    @(39, 0.08),  # AM12
    @(40, 0.284),  # AN12
    @(41, 1.118)  # AO12
)

$row13 = @(
    @(2, 3.389),  # B13
    @(3, 1.404),  # C13
    @(4, 1.185),  # D13
    @(5, 4.567),  # E13
    @(6, 0.446),  # F13
    @(7, 0.667),  # G13
    @(8, 4.588),  # H13
    @(9, 0.654),  # I13
    @(10, 0.8090000000000001),  # J13
    @(11, 2.265),  # K13
    @(12, 0.606),  # L13
    @(13, 0.779),  # M13
    @(14, 3.222),  # N13
    @(15, 0.728),  # O13
    @(16, 0.853),  # P13
    @(26, 2.514),  # Z13
    @(27, 2.878),  # AA13
    @(28, 1.697),  # AB13
    @(29, 6.4),  # AC13
    @(30, 2.24),  # AD13
    @(31, 1.497),  # AE13
    @(32, 1.556),  # AF13
    @(33, 0.58),  # AG13
    @(34, 0.762),  # AH13
    @(35, 1.194),  # AI13
    @(36, 0.157),  # AJ13
    @(37, 0.396),  # AK13
    @(38, 1.528),  # AL13
    @(39, 0.694),  # AM13
    @(41, 1.426)  # AO13
)

$rowUpdates = @(
    @(4, $row4),
    @(5, $row5),
    @(6, $row6),
    @(7, $row7),
    @(8, $row8),
    @(9, $row9),
    @(10, $row10),
    @(11, $row11),
    @(12, $row12),
    @(13, $row13)
)

foreach ($rowUpdate in $rowUpdates) {
    $rowNum = $rowUpdate[0]
    $cellPairs = $rowUpdate[1]
    foreach ($pair in $cellPairs) {
        $ws.Cells.Item($rowNum, $pair[0]).Value = $pair[1]
    }
}
